# Auto-generated edit script: updates cryptocurrency price/volume data
# (cryptos list refresh) to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (single dot or
# none) must be forced to Text format first, otherwise Excel would silently
# coerce them into numbers (losing e.g. trailing zeros) -- the source data
# keeps these as literal text strings.
$textCells = @(
    'D4',
    'D5',
    'D8',
    'D11',
    'D15',
    'D16',
    'D17',
    'D18',
    'D19',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D30',
    'D31',
    'D35',
    'D40',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D50',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell-by-cell.
$ws.Range('D2').Value = '29.145.08'
$ws.Range('E2').Value = '  +2.97%  '
$ws.Range('D3').Value = '1.580.63'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '212.44'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  +5.68%  '
$ws.Range('D8').Value = '26.35'
$ws.Range('E8').Value = '  +11.01%  '
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').Value = '0.0905'
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('D12').Value = '1.806.03'
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('D13').Value = '1.591.80'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '29.168.66'
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '3.71'
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.523'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '62.25'
$ws.Range('D18').Value = '236.53'
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').Value = '7.45'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').Value = '  +2.36%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  +2.01%  '
$ws.Range('D23').Value = '9.18'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +4.17%  '
$ws.Range('D25').Value = '153.60'
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('D26').Value = '15.15'
$ws.Range('E26').Value = '  +2.78%  '
$ws.Range('D27').Value = '0.107'
$ws.Range('E27').Value = '  +4.11%  '
$ws.Range('E28').Value = '  +1.73%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').Value = '0.0469'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').Value = '1.07'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('D33').Value = '1.422.63'
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').Value = '1.05'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('E37').Value = '  +5.91%  '
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('D40').Value = '0.531'
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('E41').Value = '  +2.38%  '
$ws.Range('B42').Value = 'BitcoinSV'
$ws.Range('C42').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D42').Value = '53.12'
$ws.Range('E42').Value = '  +25.99%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '0.790'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '0.0471'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('D46').Value = '64.72'
$ws.Range('E46').Value = '  +4.61%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').Value = '1.718.15'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('E49').Value = '  -6.61%  '
$ws.Range('D50').Value = '85.70'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  -1.25%  '
